# Updates odds values on Sheet1 (rows 2-4) to reflect the latest
# FlashScore data refresh, per commit "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 1.91
$ws.Range("H2").Value = 3
$ws.Range("J2").Value = 2.75
$ws.Range("O2").Value = 1.57
$ws.Range("P2").Value = 2.25
$ws.Range("Q2").Value = 2.88
$ws.Range("R2").Value = 1.4
$ws.Range("S2").Value = 1.62
$ws.Range("T2").Value = 2.2
$ws.Range("U2").Value = 2.5
$ws.Range("V2").Value = 1.5
$ws.Range("X2").Value = 7.5
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 17
$ws.Range("AC2").Value = 5.5
$ws.Range("AG2").Value = 8.5
$ws.Range("AN2").Value = 3.75
$ws.Range("AO2").Value = 12
$ws.Range("AT2").Value = 2.2
$ws.Range("AZ2").Value = 101

# Row 3
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75

# Row 4
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 2.5
$ws.Range("L4").Value = 6.5
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("X4").Value = 6.5
$ws.Range("AG4").Value = 10
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 67
$ws.Range("AK4").Value = 51
$ws.Range("AW4").Value = 7

$wb.Save()
